# Append four new match rows (67-70) to the Azerbaijan Premier League 2023-2024
# sheet, mirroring the layout/styling of the existing rows (row 66 is the
# last existing data row and is used as the style template).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = $ws.Range("A66:V66")

$newRows = @(
    @{
        row = 67; idx = 66
        F = "Sabail"; G = 0; H = "Kapaz"; I = 0
        E = 45255.45833333334
        J = 1.84; K = "23/11/2023 23:12"
        L = 2.4;  M = "25/11/2023 10:51"
        N = 3.24; O = "23/11/2023 23:12"
        P = 3.27; Q = "25/11/2023 10:48"
        R = 3.8;  S = "23/11/2023 23:12"
        T = 2.86; U = "25/11/2023 10:51"
        V = "https://www.betexplorer.com/football/azerbaijan/premier-league/sabail-kapaz/K8iqs7K3/"
    },
    @{
        row = 68; idx = 67
        F = "Gabala"; G = 4; H = "Turan"; I = 0
        E = 45255.54166666666
        J = 2.08; K = "24/11/2023 01:12"
        L = 2.24; M = "25/11/2023 12:09"
        N = 3.06; O = "24/11/2023 01:12"
        P = 3.18; Q = "25/11/2023 12:09"
        R = 3.25; S = "24/11/2023 01:12"
        T = 3.23; U = "25/11/2023 12:09"
        V = "https://www.betexplorer.com/football/azerbaijan/premier-league/gabala-turan/OtbzqTki/"
    },
    @{
        row = 69; idx = 68
        F = "Neftci Baku"; G = 1; H = "Araz"; I = 1
        E = 45256.52083333334
        J = 1.88; K = "25/11/2023 00:42"
        L = 2.16; M = "26/11/2023 12:26"
        N = 3.14; O = "25/11/2023 00:42"
        P = 3.22; Q = "26/11/2023 12:27"
        R = 3.75; S = "25/11/2023 00:42"
        T = 3.35; U = "26/11/2023 12:26"
        V = "https://www.betexplorer.com/football/azerbaijan/premier-league/neftci-baku-araz-pfk/WnBlLPcp/"
    },
    @{
        row = 70; idx = 69
        F = "Sumqayit"; G = 0; H = "Zira"; I = 0
        E = 45256.625
        J = 3.22; K = "25/11/2023 03:12"
        L = 3.48; M = "26/11/2023 14:55"
        N = 2.82; O = "25/11/2023 03:12"
        P = 2.79; Q = "26/11/2023 14:59"
        R = 2.23; S = "25/11/2023 03:12"
        T = 2.36; U = "26/11/2023 14:55"
        V = "https://www.betexplorer.com/football/azerbaijan/premier-league/sumqayit-fk-zira-fk/AwAhKqDj/"
    }
)

foreach ($data in $newRows) {
    $r = $data.row

    # Copy the full row (values + styles) from the last existing row so the
    # new row inherits the same cell formatting (border/bold index column,
    # date-time number format column, etc.)
    $srcRow.Copy($ws.Range("A" + $r + ":V" + $r))

    $ws.Range("A$r").Value = $data.idx
    $ws.Range("B$r").Value = "azerbaijan"
    $ws.Range("C$r").Value = "premier-league"
    $ws.Range("D$r").Value = "2023-2024"
    $ws.Range("E$r").Value = $data.E
    $ws.Range("F$r").Value = $data.F
    $ws.Range("G$r").Value = $data.G
    $ws.Range("H$r").Value = $data.H
    $ws.Range("I$r").Value = $data.I
    $ws.Range("J$r").Value = $data.J
    $ws.Range("K$r").Value = $data.K
    $ws.Range("L$r").Value = $data.L
    $ws.Range("M$r").Value = $data.M
    $ws.Range("N$r").Value = $data.N
    $ws.Range("O$r").Value = $data.O
    $ws.Range("P$r").Value = $data.P
    $ws.Range("Q$r").Value = $data.Q
    $ws.Range("R$r").Value = $data.R
    $ws.Range("S$r").Value = $data.S
    $ws.Range("T$r").Value = $data.T
    $ws.Range("U$r").Value = $data.U
    $ws.Range("V$r").Value = $data.V
}

$ws.Range("A1").Select()
